$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Orders")
$ws2 = $wb.Worksheets.Item("Daily Summary")

# A new order came in -> insert a fresh row right under the header (row 2),
# pushing the existing orders (old rows 2-7) down to rows 3-8.
$ws.Rows.Item(2).Insert()

# Fill in the new order's data.
$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = "2026-01-13 16:54"
$ws.Cells.Item(2, 3).Value = "Pooja"
$ws.Cells.Item(2, 4).Value = "a14"

# Phone numbers are all-digit strings, so Excel would otherwise coerce the
# value to a Number - force it to stay text like the rest of the column,
# then drop back to the Normal style so no stray formatting is left behind.
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "9096648553"
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(2, 6).Value = "Wheat Chapati x1"
$ws.Cells.Item(2, 7).Value = 15
$ws.Cells.Item(2, 8).Value = "NEW"
$ws.Cells.Item(2, 9).Value = "PENDING"

# Collection Date/Time aren't set yet for this brand-new order, and
# Notes/Cancel Reason/Feedback stay blank too - but keep them as empty
# text cells (matching every other row in this sheet) rather than
# leaving them completely untouched/empty.
$emptyTextCols = 10,11,12,13,14
foreach ($col in $emptyTextCols) {
    $ws.Cells.Item(2, $col).Value = "'"
    $ws.Cells.Item(2, $col).Style = "Normal"
}

# Update today's Daily Summary roll-up for the new order.
$ws2.Cells.Item(2, 2).Value = 7
$ws2.Cells.Item(2, 5).Value = 195
$ws2.Cells.Item(2, 7).Value = 195
